# One more run after removing Brown-crowned Scimitar-Babbler
# Decrement the affected summary counts by 1 across the summary sheets.

$wb = $excel.ActiveWorkbook

# Sheet: Trends Status -> row "eBird Data Deficient"
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B4").Value = 422
$ws1.Range("C4").Value = 302

# Sheet: Range Status -> row "Very Restricted"
$ws2 = $wb.Worksheets.Item("Range Status")
$ws2.Range("B3").Value = 69

# Sheet: Priority Status -> row "High"
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 198

# Sheet: Species qualification -> "SoIB 2023 Assessment" and "Range Analysis"
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B2").Value = 945
$ws4.Range("B5").Value = 945

# Sheet: SoIB-IUCN cross-tab -> row "Least Concern" and row "Sum"
$ws5 = $wb.Worksheets.Item("SoIB-IUCN cross-tab")
$ws5.Range("B6").Value = 100
$ws5.Range("E6").Value = 786
$ws5.Range("B8").Value = 198
$ws5.Range("E8").Value = 945
